$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row at position 11 for "Jurisdiction", pushing Description and
# everything below it down by one row.
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Update the Version value (now row 3)
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# Update the Date value (now row 8)
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"
